$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Expand the "API Calls" paragraph: replace the old sentence about the AI
#    coach with the new text describing front-end/back-end API routes.
# ---------------------------------------------------------------------------
$oldText = " enable our AI coach (more on that in Features). We will also make use of API calls "
$newText = " enable our front-end features to connect to our application’s back-end. Fetching user and activity information stored in our database will be possible by creating multiple API routes. This offers both clarity and organization for development and usability. API calls will be used "
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Add a trailing space after the "...etc." sentence at the end of that
#    same paragraph, then split a new, empty paragraph in right after it.
# ---------------------------------------------------------------------------
$tail = $d.Content
$tail.Find.Execute("display the information that the user wants to see based on filters, searches, organizations, etc.") | Out-Null
$tail.Collapse(0)
$tail.InsertAfter(" ")
$tail.Collapse(0)
$tail.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 3. Move the page break: it now renders before the "VS Code" heading
#    instead of in the middle of "Vs Code's vast extensions".
# ---------------------------------------------------------------------------
$heading = $d.Content
$heading.Find.Execute("VS Code") | Out-Null
$heading.Collapse(1)
$heading.InsertBefore([char]2)
$markerRange = $d.Content
$markerRange.Find.Execute([char]2) | Out-Null
$markerRange.Delete() | Out-Null

$vsCode = $d.Content
$vsCode.Find.Execute("Vs Code’s vast extensions", $true) | Out-Null
$vsCode.Text = "Vs Code’s vast extensions"
